# Daily attendance processing - 2026-01-19 21:02:59
# Swap the order of the "Recorded By" entries in column G from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# wherever they occur on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
